# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Recalculated K values replace the previous (stale) Strike# numbers for
# rows 2-30 (the 29 most recent outings), column G ("K" header in row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(5, 6, 2, 6, 6, 4, 9, 2, 4, 2, 2, 0, 0, 1, 4, 2, 6, 4, 2, 4, 3, 7, 6, 3, 7, 5, 3, 4, 4)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
